# Applies the cryptos-list refresh described in the commit:
# "Updated cryptos list on Thu Aug 22 03:50:47 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.176.29"
$ws.Range("E2").Value = "  +1.70%  "
$ws.Range("D3").Value = "2.588.65"
$ws.Range("E3").Value = "  +0.27%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'567.62"
$ws.Range("E5").Value = "  +0.74%  "
$ws.Range("D6").Value = "'141.28"
$ws.Range("E6").Value = "  -0.80%  "
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").Value = "2.607.20"
$ws.Range("E9").Value = "  +0.72%  "
$ws.Range("D10").Value = "'6.58"
$ws.Range("E10").Value = "  -0.86%  "
$ws.Range("E11").Value = "  +1.22%  "
$ws.Range("E12").Value = "  +5.98%  "
$ws.Range("E13").Value = "  -6.14%  "
$ws.Range("D14").Value = "3.050.13"
$ws.Range("E14").Value = "  +0.42%  "
$ws.Range("D15").Value = "60.191.98"
$ws.Range("E15").Value = "  +1.78%  "
$ws.Range("D16").Value = "'23.20"
$ws.Range("E16").Value = "  +0.93%  "
$ws.Range("E17").Value = "  +1.92%  "
$ws.Range("D18").Value = "2.598.86"
$ws.Range("E18").Value = "  +0.61%  "
$ws.Range("D19").Value = "'11.34"
$ws.Range("E19").Value = "  +9.29%  "
$ws.Range("E20").Value = "  +1.85%  "
$ws.Range("D21").Value = "'345.51"
$ws.Range("E21").Value = "  +2.54%  "
$ws.Range("D22").Value = "'6.91"
$ws.Range("E22").Value = "  +8.54%  "
$ws.Range("D24").Value = "'0.538"
$ws.Range("E24").Value = "  +16.29%  "
$ws.Range("D25").Value = "'62.97"
$ws.Range("E25").Value = "  -1.80%  "
$ws.Range("D26").Value = "'0.997"
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("E27").Value = "  -2.13%  "
$ws.Range("D28").Value = "'7.64"
$ws.Range("E28").Value = "  +4.69%  "
$ws.Range("D29").Value = "0.0$([char]0x2083)0779"
$ws.Range("E29").Value = "  +0.69%  "
$ws.Range("E30").Value = "  +7.21%  "
$ws.Range("E31").Value = "  -0.09%  "
$ws.Range("D32").Value = "'6.30"
$ws.Range("E32").Value = "  +3.58%  "
$ws.Range("D33").Value = "'160.38"
$ws.Range("E33").Value = "  -0.75%  "
$ws.Range("D34").Value = "'19.42"
$ws.Range("E34").Value = "  +2.52%  "
$ws.Range("E35").Value = "  +4.78%  "
$ws.Range("D36").Value = "'0.958"
$ws.Range("E36").Value = "  +9.46%  "
$ws.Range("E37").Value = "  +4.35%  "
$ws.Range("E38").Value = "  +7.51%  "
$ws.Range("E39").Value = "  +0.79%  "
$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D40").Value = "'3.80"
$ws.Range("E40").Value = "  +3.90%  "
$ws.Range("B41").Value = "SuiNetwork"
$ws.Range("C41").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D41").Value = "'0.849"
$ws.Range("E41").Value = "  -2.52%  "
$ws.Range("D42").Value = "'293.59"
$ws.Range("E42").Value = "  -0.33%  "
$ws.Range("D43").Value = "'138.51"
$ws.Range("E43").Value = "  +4.86%  "
$ws.Range("E45").Value = "  +1.43%  "
$ws.Range("D46").Value = "'0.0976"
$ws.Range("E46").Value = "  +0.29%  "
$ws.Range("D47").Value = "'19.58"
$ws.Range("E47").Value = "  +2.89%  "
$ws.Range("D48").Value = "'0.0544"
$ws.Range("E48").Value = "  +1.73%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "'0.0239"
$ws.Range("E49").Value = "  +2.58%  "
$ws.Range("B50").Value = "WhiteBITCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D50").Value = "'10.67"
$ws.Range("E50").Value = "  +0.25%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").Value = "'19.60"
$ws.Range("E51").Value = "  +5.78%  "
